$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "88÷6="
$t.Cell(1, 2).Range.Text = "39÷7="
$t.Cell(1, 3).Range.Text = "41÷4="
$t.Cell(1, 4).Range.Text = "22÷6="
$t.Cell(1, 5).Range.Text = "24÷3="
$t.Cell(5, 1).Range.Text = "35÷4="
$t.Cell(5, 2).Range.Text = "75÷6="
$t.Cell(5, 3).Range.Text = "78÷7="
$t.Cell(5, 4).Range.Text = "29÷3="
$t.Cell(5, 5).Range.Text = "83÷2="
$t.Cell(9, 1).Range.Text = "93÷4="
$t.Cell(9, 2).Range.Text = "21÷5="
$t.Cell(9, 3).Range.Text = "37÷8="
$t.Cell(9, 4).Range.Text = "73÷6="
$t.Cell(9, 5).Range.Text = "18÷7="
$t.Cell(13, 1).Range.Text = "23÷8="
$t.Cell(13, 2).Range.Text = "14÷3="
$t.Cell(13, 3).Range.Text = "85÷8="
$t.Cell(13, 4).Range.Text = "68÷8="
$t.Cell(13, 5).Range.Text = "98÷6="
$t.Cell(17, 1).Range.Text = "22÷6="
$t.Cell(17, 2).Range.Text = "58÷7="
$t.Cell(17, 3).Range.Text = "94÷5="
$t.Cell(17, 4).Range.Text = "37÷7="
$t.Cell(17, 5).Range.Text = "79÷7="
